$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 5641
$wsExpo.Range("F5").Value  = 7
$wsExpo.Range("F8").Value  = 2533
$wsExpo.Range("F9").Value  = 84
$wsExpo.Range("F10").Value = 153
$wsExpo.Range("F11").Value = 12
$wsExpo.Range("F12").Value = 80
$wsExpo.Range("F13").Value = 21
$wsExpo.Range("F14").Value = 2366
$wsExpo.Range("F15").Value = 364

# Sheet "全部类型" (all types) - same underlying events, different row numbers
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 5641
$wsAll.Range("F6").Value  = 7
$wsAll.Range("F10").Value = 2533
$wsAll.Range("F11").Value = 84
$wsAll.Range("F12").Value = 153
$wsAll.Range("F13").Value = 12
$wsAll.Range("F15").Value = 80
$wsAll.Range("F16").Value = 21
$wsAll.Range("F17").Value = 2366
$wsAll.Range("F18").Value = 364

$wb.Save()
